# Update automatico via Actualizar 03-11-2021 12-01-22
# Shifts the "Fecha" (Ultimo/last-checked) timestamps in column D down one
# block (14 rows per availability check) and stamps the newest block with
# the latest check time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value  = 44266.5006951052
$ws.Range("D16:D29").Value = 44266.47930578703
$ws.Range("D30:D43").Value = 44266.4579140625
